# Add data for 2022-09-18
# - Advances the "through" date from 2022-09-09 to 2022-09-10 (sheet name +
#   header label), and bumps the per-neighborhood / per-month carjacking
#   counts for the newly-included day across the current month (column B)
#   and the matching "September <year>" columns for prior years.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet title / running header -----------------------------------------
$wb.Worksheets.Item(1).Name = "Through 2022-09-10"
$ws.Range("B1").Value = "September 2022 (through September 10)"

# --- Helper to bump an existing numeric cell by +1 -------------------------
function Add-One([string]$addr) {
    $cell = $ws.Range($addr)
    $cell.Value = $cell.Value2 + 1
}

# Row 2 - Austin
Add-One "K2"
Add-One "AL2"

# Row 3 - Garfield Park
Add-One "B3"

# Row 4 - Auburn Gresham
$ws.Range("B4").Value = 1
$ws.Range("T4").Value = 1

# Row 8 - Humboldt Park
Add-One "B8"
Add-One "K8"
Add-One "T8"
$ws.Range("AU8").Value = 1

# Row 10 - North Lawndale
Add-One "K10"

# Row 11 - West Pullman
Add-One "AU11"

# Row 12 - Grand Crossing
Add-One "K12"

# Row 14 - Roseland
Add-One "B14"

# Row 15 - Chatham
$ws.Range("B15").Value = 1
Add-One "AC15"

# Row 25 - Calumet Heights
$ws.Range("B25").Value = 1

# Row 29 - West Loop
Add-One "B29"
Add-One "T29"

# Row 30 - Bridgeport
$ws.Range("BD30").Value = 1

# Row 37 - Rogers Park
$ws.Range("K37").Value = 1

# Row 38 - Wicker Park
$ws.Range("K38").Value = 1

# Row 50 - Grand Boulevard
Add-One "AU50"

# Row 64 - Douglas
Add-One "B64"

# Row 65 - Dunning
$ws.Range("T65").Value = 1

# Row 66 - East Village
$ws.Range("AC66").Value = 1

# Row 91 - Portage Park
$ws.Range("BM91").Value = 1

# Row 96 - Washington Heights
$ws.Range("BD96").Value = 1

# Row 97 - Washington Park
$ws.Range("B97").Value = 1
